$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.695.26'
$ws.Range('E2').Value = '  +0.41%  '

$ws.Range('D3').Value = '3.377.72'
$ws.Range('E3').Value = '  -0.67%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.04%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '175.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.17%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('D8').Value = '3.374.29'
$ws.Range('E8').Value = '  -0.57%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.06%  '

$ws.Range('E10').Value = '  +1.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.632'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.58%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.37%  '

$ws.Range('D15').Value = '3.924.90'
$ws.Range('E15').Value = '  -0.60%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.403.48'
$ws.Range('E16').Value = '  -0.02%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.18'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.90%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.119'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.18%  '

$ws.Range('D19').Value = '65.706.31'
$ws.Range('E19').Value = '  +0.35%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '466.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.40%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.52'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.11%  '

$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.91%  '

$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.56%  '

$ws.Range('E27').Value = '  +0.30%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.19%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.94%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.43%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.58'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.42'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.17%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '577.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '62.02'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.79%  '

$ws.Range('E35').Value = '  -0.83%  '

$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('E37').Value = '  +2.03%  '

$ws.Range('E38').Value = '  +0.87%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.03%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.20%  '

$ws.Range('D41').Value = '0.0₃0738'
$ws.Range('E41').Value = '  -2.53%  '

$ws.Range('D42').Value = '3.095.05'
$ws.Range('E42').Value = '  +0.08%  '

$ws.Range('E43').Value = '  -1.06%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0416'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.53%  '

$ws.Range('E45').Value = '  -1.24%  '

$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.90%  '

$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.20%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.19%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.39'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.64%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.88%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.76%  '
